$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iOS_Estimate")

$dt = Get-Date -Year 2015 -Month 7 -Day 30 -Hour 0 -Minute 0 -Second 0

$ws.Range("F18").NumberFormat = "d-mmm"
$ws.Range("F18").Value = $dt
$ws.Range("H18").Value = "In progress"

$ws.Range("F19").NumberFormat = "d-mmm"
$ws.Range("F19").Value = $dt
$ws.Range("H19").Value = "In progress"

$ws.Range("F34").NumberFormat = "d-mmm"
$ws.Range("F34").Value = $dt
$ws.Range("G34").NumberFormat = "d-mmm"
$ws.Range("G34").Value = $dt
$ws.Range("H34").Value = "completed"

$ws.Range("H20").Select()
